# This workbook's sheet is protected, so it must be unprotected before any
# cell can be edited, and re-protected afterwards to preserve the original
# (protected) state of the sheet.
$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer note (cell A10)
$oldText = $ws.Range("A10").Value2
$newText = $oldText -replace "2021-03-24", "2021-03-25"
$ws.Range("A10").Value = $newText

# Update the Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.4939899275103433
$ws.Range("E2").Value = -0.001576044129235621

$ws.Range("D3").Value = 0.3289040501242735
$ws.Range("E3").Value = 0.01224820325943932

$ws.Range("D4").Value = 0.09332867227786038
$ws.Range("E4").Value = 0.003457548981943948

$ws.Range("D5").Value = 0.05538844903379138
$ws.Range("E5").Value = 0.001733903594959951

$ws.Range("D6").Value = 0.02838890105373147
$ws.Range("E6").Value = 0.007442489851150258

$ws.Range("E7").Value = 0.003879944528317969

# Re-protect the sheet (password hashes to the same legacy "D382" value as
# the original workbook's protection).
$ws.Protect("ZYh3qiPB1")
